$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 10: "动态模糊" (Motion Blur) feature entry, matching the
# "not yet implemented" formatting used by rows like B3/B5/B6/B8.
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "动态模糊"
$ws.Range("E10").Value = "主要是录屏需要，帧率低的时候好看"

# Clone the formatting (style) from an existing "todo" row so B10 matches
# the same visual style (font/color) used for other pending items.
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "动态模糊"

# Update the active selection to reflect where the user left off editing.
[void]$ws.Range("D11").Select()
